$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Correct the two existing tail rows (325, 326) ---
$ws.Range("C325").Value = 7710432000000
$ws.Range("D325").Value = 7710432000000
$ws.Range("E325").Value = 7710432000000
$ws.Range("F325").Value = 7710432000000

$ws.Range("C326").Value = 7822810000000
$ws.Range("D326").Value = 7822810000000
$ws.Range("E326").Value = 7822810000000
$ws.Range("F326").Value = 7822810000000

# --- Append three new monthly rows (327-329) ---
$newRows = @(
    @{ Row = 327; DateSerial = 44986.45833333334; Value = 7965088000000 },
    @{ Row = 328; DateSerial = 45017.45833333334; Value = 8069151000000 },
    @{ Row = 329; DateSerial = 45047.41666666666; Value = 8140535000000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting from the row above (keeps the date number format,
    # font, border, alignment consistent with the rest of the column).
    $ws.Range("A" + ($row - 1)).Copy() | Out-Null
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $aCell.Value = $r.DateSerial

    $ws.Cells.Item($row, 2).Value = "ECONOMICS:EGM2"
    $ws.Cells.Item($row, 3).Value = $r.Value
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}

$excel.CutCopyMode = $false
